# Apply the change described by the diff:
#  - Add a new project row (row 9) to Sheet1 with mentor Arun's
#    "Anomaly Detection using Adversarial Learning on Big Data" project.
#  - Adjust row heights for rows 6, 8 and the new row 9.
#  - Update sheet view selection (D3) and page orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Copy formatting from the last existing data row (row 8) down into
#    the new row 9, so fonts/alignment/number formats match the rest of
#    the table.
# ---------------------------------------------------------------------
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)   # xlPasteFormats

# Column A keeps the plain (unformatted / column-default) style rather
# than the "Sr. No" style used on other rows, so pull that look from an
# untouched cell further down column A (still using the column default).
$ws.Range("A20").Copy()
$ws.Range("A9").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 152
$ws.Rows.Item(8).RowHeight = 171
$ws.Rows.Item(9).RowHeight = 409.6

# ---------------------------------------------------------------------
# 3. Cell values for the new row (9)
#    NOTE: columns E/G are populated before F so that the new shared
#    strings are interned in the same order as the source workbook.
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Arun"
$ws.Cells.Item(9, 3).Value = "Anomaly Detection using Adversarial Learning on Big Data"
$ws.Cells.Item(9, 4).Value = "Identification of events from streaming videos, audio or other big data, which does not conform with the training data distribution and expected behavior is essential for modern smart cities. As we embrace inclusion of security cameras, digital assistants which listen to us, etc., the use of algorithms which detect anomalies could help deter threats, and unwanted scenarios automatically. In industries, automated anomaly detection algorithms could learn from a set of training data to monitor specific tasks, altering in cases of anomalies in different conditions – that might not be possible with present control systems. For example, proper working of a product line, or gear mechanisms could be monitored using anomaly detection algorithms. This project focuses on a Generative Adversarial Network (GAN) based anomaly detection algorithm which learns to generalize on a training data distribution and expected behavior of objects in the dataset. Hence, once trained, the GAN model can discriminate anomalies and alert users on the onset of anomalies. The main goal of the project is to use a GAN architecture to find anomalies in the data sources mentioned below."
$ws.Cells.Item(9, 5).Value = "Data Sources:`nCUHK Avenue dataset`nCUHK Avenue dataset contains 16 training videos and 21 testing ones with a total of 47 abnormal events, including throwing objects, loitering and running. The size of people may change because of the camera position and angle`nUCSD Pedestrian dataset`nThe UCSD dataset contains two parts: The UCSD Pedestrian 1 (Ped1) dataset and the UCSD Pedestrian 2 (Ped2) dataset. The UCSD Pedestrian 1 (Ped1) dataset includes 34 training videos and 36 testing ones with 40 irregular events. All of these abnormal cases are about vehicles such as bicycles and cars. The UCSD Pedestrian 2 (Ped2) dataset contains 16 training videos and 12 testing videos with 12 abnormal events. The definition of anomaly for Ped2 is the same with Ped1. Usually different methods are evaluated on these two parts separately`nShanghaiTech dataset`nThe ShanghaiTech dataset is a very challenging anomaly detection dataset. It contains 330 training videos and 107 testing ones with 130 abnormal events. Totally, it consists of 13 scenes and various anomaly types."
$ws.Cells.Item(9, 7).Value = "Not required."
$ws.Cells.Item(9, 6).Value = "Tasks:`nDownload the datasets from https://github.com/StevenLiuWen/ano_pred_cvpr2018#2-download-datasets`nRead and understand about anomaly detection using GANs: https://arxiv.org/pdf/1712.09867.pdf`nWork with Arun and run the training code available in: https://github.com/StevenLiuWen/ano_pred_cvpr2018#4-training-from-scratch-here-we-use-ped2-and-avenue-datasets-for-examples`nCollect your own dataset and test the model – fun experiments! Learn, and explore scenarios to apply the algorithm from a bigdata context.`nWrite a short report on your findings, and experiments."

# ---------------------------------------------------------------------
# 4. Sheet view: select D3 (this also clears the old "topLeftCell=A5"
#    scroll position that was saved with the workbook).
# ---------------------------------------------------------------------
[void]$ws.Range("D3").Select()

# ---------------------------------------------------------------------
# 5. Page orientation -> portrait
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
